# The sheet is a daily/weekly log of Caqui (persimmon) price observations at
# Macroferia Regional de Talca, sorted newest-first within the dataset.
# A new weekly observation (2022-06-09) was added at the top of the data
# block (row 23), pushing every existing row down by one and extending the
# used range from A1:T50 to A1:T51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 - this shifts rows 23:50 down to 24:51
# and automatically extends the sheet dimension to A1:T51.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(23, 1).Value = 5
$ws.Cells.Item(23, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(23, 3).Value = "Maule"
$ws.Cells.Item(23, 4).Value = 44721
$ws.Cells.Item(23, 5).Value = 7
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100107
$ws.Cells.Item(23, 8).Value = "Otros"
$ws.Cells.Item(23, 9).Value = 100107001
$ws.Cells.Item(23, 10).Value = "Caqui"
$ws.Cells.Item(23, 11).Value = "Mankaki"
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 80
$ws.Cells.Item(23, 14).Value = 12000
$ws.Cells.Item(23, 15).Value = 12000
$ws.Cells.Item(23, 16).Value = 12000
$ws.Cells.Item(23, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(23, 18).Value = "Región del Maule"
$ws.Cells.Item(23, 19).Value = 667
$ws.Cells.Item(23, 20).Value = 18
